$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain decimal numbers need to be forced to Text
# format before assignment (otherwise Excel auto-converts "561.00" -> 561 etc.),
# then the format is reset back to Normal/General so no style/format diff remains.
$textCells = @("D5", "D6", "D11", "D13", "D14", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D37", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.135.35"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "2.410.10"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "561.00"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "138.46"
$ws.Range("E6").Value = "  +4.85%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "2.407.66"
$ws.Range("E9").Value = "  +3.09%  "
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").Value = "5.72"
$ws.Range("E11").Value = "  +3.59%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("D14").Value = "25.85"
$ws.Range("E14").Value = "  +8.23%  "
$ws.Range("D15").Value = "2.835.99"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").Value = "62.063.88"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("E17").Value = "  +3.89%  "
$ws.Range("D18").Value = "2.412.85"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").Value = "11.05"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").Value = "343.85"
$ws.Range("E20").Value = "  +8.84%  "
$ws.Range("D21").Value = "4.24"
$ws.Range("E21").Value = "  +1.53%  "
$ws.Range("D22").Value = "6.89"
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "65.17"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "8.38"
$ws.Range("E27").Value = "  +6.15%  "
$ws.Range("D28").Value = "1.52"
$ws.Range("E28").Value = "  +11.03%  "
$ws.Range("D29").Value = "1.37"
$ws.Range("E29").Value = "  +13.52%  "
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("D31").Value = "0.0₃0774"
$ws.Range("E31").Value = "  +4.62%  "
$ws.Range("D32").Value = "6.39"
$ws.Range("E32").Value = "  +7.14%  "
$ws.Range("D33").Value = "171.68"
$ws.Range("D34").Value = "1.42"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D36").Value = "18.56"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "4.54"
$ws.Range("E37").Value = "  +10.47%  "
$ws.Range("D38").Value = "366.31"
$ws.Range("E38").Value = "  +10.89%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +7.42%  "
$ws.Range("D42").Value = "39.02"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").Value = "144.27"
$ws.Range("E43").Value = "  +2.90%  "
$ws.Range("D44").Value = "3.68"
$ws.Range("D45").Value = "20.65"
$ws.Range("E45").Value = "  +6.38%  "
$ws.Range("D46").Value = "0.0967"
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").Value = "0.0519"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").Value = "0.584"
$ws.Range("E48").Value = "  +3.63%  "
$ws.Range("E49").Value = "  +3.55%  "
$ws.Range("D50").Value = "17.93"
$ws.Range("E50").Value = "  +5.02%  "
$ws.Range("D51").Value = "0.0₆0216"
$ws.Range("E51").Value = "  -3.34%  "

foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}
